$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.228.72'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '1.773.16'
$ws.Range('E3').Value = '  +3.51%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.85'
$ws.Range('E5').Value = '  +1.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5259'
$ws.Range('E7').Value = '  +10.79%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3677'
$ws.Range('E8').Value = '  +6.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.71'
$ws.Range('E9').Value = '  +1.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07353'
$ws.Range('E10').Value = '  +0.95%  '
$ws.Range('E11').Value = '  +4.16%  '
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('E13').Value = '  +2.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.058'
$ws.Range('E14').Value = '  +3.32%  '
$ws.Range('D15').Value = '1.766.07'
$ws.Range('E15').Value = '  +3.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.926'
$ws.Range('E16').Value = '  +1.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.69'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('E18').Value = '  +0.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06437'
$ws.Range('E19').Value = '  +1.23%  '
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.70'
$ws.Range('E21').Value = '  +1.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.805'
$ws.Range('E22').Value = '  +3.58%  '
$ws.Range('D23').Value = '27.276.71'
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.26'
$ws.Range('E24').Value = '  +4.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.108'
$ws.Range('E25').Value = '  +0.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.35'
$ws.Range('E26').Value = '  +1.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.10'
$ws.Range('E27').Value = '  +1.93%  '
$ws.Range('D28').Value = '1.973.42'
$ws.Range('E28').Value = '  +3.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.326'
$ws.Range('E29').Value = '  +11.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '121.03'
$ws.Range('E30').Value = '  +0.76%  '
$ws.Range('E31').Value = '  +4.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09776'
$ws.Range('E32').Value = '  +5.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.558'
$ws.Range('E33').Value = '  +4.76%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.622'
$ws.Range('E34').Value = '  +0.92%  '
$ws.Range('E35').Value = '  +1.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05964'
$ws.Range('E36').Value = '  +1.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '11.20'
$ws.Range('E37').Value = '  +1.28%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2021'
$ws.Range('E38').Value = '  +0.34%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.828'
$ws.Range('E39').Value = '  +1.68%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6127'
$ws.Range('E40').Value = '  +3.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.431'
$ws.Range('E41').Value = '  +1.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.047'
$ws.Range('E42').Value = '  +7.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.138'
$ws.Range('E43').Value = '  +2.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.12'
$ws.Range('E44').Value = '  +3.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5750'
$ws.Range('E45').Value = '  +2.13%  '
$ws.Range('E46').Value = '  +1.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '121.24'
$ws.Range('E47').Value = '  +2.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.875'
$ws.Range('E48').Value = '  +1.64%  '
$ws.Range('E49').Value = '  +2.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06692'
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('E51').Value = '  +0.05%  '
